$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (not numeric) for numeric-looking Price values,
# matching the source file where these are plain text cells.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values from the crypto-price refresh.
$ws.Range("D2").Value = "60.441.32"
$ws.Range("E2").Value = "  -1.73%  "
$ws.Range("D3").Value = "2.334.69"
$ws.Range("E3").Value = "  -4.38%  "
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").Value = "541.45"
$ws.Range("E5").Value = "  -1.03%  "
$ws.Range("D6").Value = "135.45"
$ws.Range("E6").Value = "  -7.06%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("D8").Value = "0.519"
$ws.Range("E8").Value = "  -10.94%  "
$ws.Range("D9").Value = "2.333.38"
$ws.Range("E9").Value = "  -4.38%  "
$ws.Range("E10").Value = "  -2.34%  "
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("E12").Value = "  -3.38%  "
$ws.Range("E13").Value = "  -3.04%  "
$ws.Range("D14").Value = "24.31"
$ws.Range("E14").Value = "  -5.88%  "
$ws.Range("E15").Value = "  -4.55%  "
$ws.Range("D16").Value = "60.166.40"
$ws.Range("E17").Value = "  -5.40%  "
$ws.Range("D18").Value = "2.329.34"
$ws.Range("E18").Value = "  -4.55%  "
$ws.Range("D19").Value = "10.52"
$ws.Range("E19").Value = "  -2.12%  "
$ws.Range("D20").Value = "315.91"
$ws.Range("E20").Value = "  -0.85%  "
$ws.Range("E21").Value = "  -1.72%  "
$ws.Range("D22").Value = "6.48"
$ws.Range("E22").Value = "  -5.71%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "62.61"
$ws.Range("E24").Value = "  -1.79%  "
$ws.Range("D25").Value = "1.65"
$ws.Range("E25").Value = "  -10.20%  "
$ws.Range("D26").Value = "8.27"
$ws.Range("E26").Value = "  +6.04%  "
$ws.Range("E27").Value = "  -0.22%  "
$ws.Range("D28").Value = "2.439.83"
$ws.Range("E28").Value = "  -4.81%  "
$ws.Range("E29").Value = "  -4.21%  "
$ws.Range("D30").Value = "1.35"
$ws.Range("E30").Value = "  -8.26%  "
$ws.Range("D31").Value = "488.78"
$ws.Range("E31").Value = "  -6.76%  "
$ws.Range("D32").Value = "0.0₃0848"
$ws.Range("E32").Value = "  -12.51%  "
$ws.Range("E33").Value = "  -2.28%  "
$ws.Range("D34").Value = "1.77"
$ws.Range("E34").Value = "  -5.12%  "
$ws.Range("E35").Value = "  -6.67%  "
$ws.Range("D36").Value = "0.995"
$ws.Range("E36").Value = "  -0.35%  "
$ws.Range("E37").Value = "  -4.35%  "
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").Value = "18.40"
$ws.Range("E38").Value = "  +1.13%  "
$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").Value = "0.371"
$ws.Range("E39").Value = "  -2.64%  "
$ws.Range("E40").Value = "  -8.83%  "
$ws.Range("E41").Value = "  +0.78%  "
$ws.Range("D42").Value = "141.16"
$ws.Range("E42").Value = "  +1.67%  "
$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  -0.22%  "
$ws.Range("D44").Value = "40.38"
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").Value = "140.17"
$ws.Range("E45").Value = "  -1.62%  "
$ws.Range("E46").Value = "  -2.40%  "
$ws.Range("D47").Value = "2.02"
$ws.Range("E47").Value = "  -10.93%  "
$ws.Range("D48").Value = "0.0507"
$ws.Range("E48").Value = "  -3.50%  "
$ws.Range("D49").Value = "18.87"
$ws.Range("E49").Value = "  -10.80%  "
$ws.Range("D50").Value = "0.564"
$ws.Range("E50").Value = "  -3.82%  "
$ws.Range("D51").Value = "0.0892"
$ws.Range("E51").Value = "  -4.04%  "

# Strip the temporary Text number-format so styling matches the original
# (cells keep their text content but no longer carry an explicit style).
$ws.Range("D4").ClearFormats()
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D7").ClearFormats()
$ws.Range("D8").ClearFormats()
$ws.Range("D14").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D26").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D36").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("D48").ClearFormats()
$ws.Range("D49").ClearFormats()
$ws.Range("D50").ClearFormats()
$ws.Range("D51").ClearFormats()
